$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.331.32'
$ws.Range("E2").Value = '  -0.82%  '
$ws.Range("D3").Value = '1.560.86'
$ws.Range("E3").Value = '  -0.88%  '
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = "'1.000"
$ws.Range("E5").Value = '  -0.15%  '
$ws.Range("D6").Value = "'287.33"
$ws.Range("E6").Value = '  +0.11%  '
$ws.Range("D7").Value = "'0.3792"
$ws.Range("E7").Value = '  +3.36%  '
$ws.Range("D8").Value = "'0.3276"
$ws.Range("E8").Value = '  -2.03%  '
$ws.Range("D9").Value = "'44.24"
$ws.Range("E9").Value = '  -9.02%  '
$ws.Range("D10").Value = "'1.149"
$ws.Range("E10").Value = '  +1.37%  '
$ws.Range("D11").Value = "'0.07390"
$ws.Range("E11").Value = '  -0.83%  '
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = '  -0.08%  '
$ws.Range("D13").Value = "'20.50"
$ws.Range("E13").Value = '  -1.78%  '
$ws.Range("D14").Value = "'5.847"
$ws.Range("E14").Value = '  -2.54%  '
$ws.Range("D15").Value = "'6.812"
$ws.Range("E15").Value = '  -1.76%  '
$ws.Range("D16").Value = '1.565.99'
$ws.Range("E16").Value = '  -0.82%  '
$ws.Range("D17").Value = "'0.00001089"
$ws.Range("E17").Value = '  -2.17%  '
$ws.Range("D18").Value = "'0.06712"
$ws.Range("E18").Value = '  -0.72%  '
$ws.Range("D19").Value = "'86.14"
$ws.Range("E19").Value = '  -2.57%  '
$ws.Range("D20").Value = "'6.396"
$ws.Range("E20").Value = '  -0.15%  '
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = '  -0.16%  '
$ws.Range("D22").Value = "'16.22"
$ws.Range("E22").Value = '  -1.59%  '
$ws.Range("D23").Value = "'11.70"
$ws.Range("E23").Value = '  -3.84%  '
$ws.Range("D24").Value = '22.340.30'
$ws.Range("E24").Value = '  -0.78%  '
$ws.Range("D25").Value = "'2.291"
$ws.Range("E25").Value = '  -3.97%  '
$ws.Range("D26").Value = "'2.588"
$ws.Range("E26").Value = '  -1.40%  '
$ws.Range("D27").Value = "'150.28"
$ws.Range("E27").Value = '  -1.52%  '
$ws.Range("D28").Value = "'19.38"
$ws.Range("E28").Value = '  -1.16%  '
$ws.Range("D29").Value = "'4.935"
$ws.Range("E29").Value = '  -1.45%  '
$ws.Range("D30").Value = "'122.48"
$ws.Range("E30").Value = '  -1.31%  '
$ws.Range("D31").Value = '1.740.66'
$ws.Range("E31").Value = '  -0.84%  '
$ws.Range("D32").Value = "'1.077"
$ws.Range("E32").Value = '  +1.91%  '
$ws.Range("D33").Value = "'5.968"
$ws.Range("E33").Value = '  -3.47%  '
$ws.Range("D34").Value = "'1.901"
$ws.Range("E34").Value = '  -5.08%  '
$ws.Range("D35").Value = "'9.547"
$ws.Range("E35").Value = '  -2.62%  '
$ws.Range("D36").Value = "'0.08282"
$ws.Range("E36").Value = '  -0.07%  '
$ws.Range("D37").Value = "'0.02392"
$ws.Range("E37").Value = '  -2.48%  '
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").Value = "'1.283"
$ws.Range("E38").Value = '  -1.25%  '
$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").Value = "'5.325"
$ws.Range("E39").Value = '  -2.46%  '
$ws.Range("D40").Value = "'0.2182"
$ws.Range("E40").Value = '  -3.93%  '
$ws.Range("D41").Value = "'0.06271"
$ws.Range("E41").Value = '  -2.90%  '
$ws.Range("D42").Value = "'11.13"
$ws.Range("E42").Value = '  -2.00%  '
$ws.Range("D43").Value = "'0.6090"
$ws.Range("E43").Value = '  -4.19%  '
$ws.Range("E44").Value = '  -0.16%  '
$ws.Range("D45").Value = "'13.93"
$ws.Range("E45").Value = '  +0.12%  '
$ws.Range("D46").Value = "'0.5935"
$ws.Range("E46").Value = '  -4.17%  '
$ws.Range("D47").Value = "'3.743"
$ws.Range("E47").Value = '  -0.46%  '
$ws.Range("D48").Value = "'2.006"
$ws.Range("E48").Value = '  -2.59%  '
$ws.Range("D49").Value = "'123.76"
$ws.Range("E49").Value = '  -1.15%  '
$ws.Range("D50").Value = "'1.178"
$ws.Range("E50").Value = '  -3.37%  '
$ws.Range("D51").Value = "'0.07107"
$ws.Range("E51").Value = '  -1.97%  '
